$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.022.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.794.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "359.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.33%  "
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  +1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.234.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.788.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.945"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.963.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  -2.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +15.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0473"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.63%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.03%  "
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0848"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.21%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.17%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -3.73%  "
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -1.36%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.09%  "
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "119.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.090.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  -4.56%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.57%  "
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.87%  "

Write-Host "Applied cryptos list update"
